# AutoCommit_21 июня 2024 г. 22:31:27_SibNout2023
#
# Adds a "ТК_оригинал" / "на момент выгрузки в элжуре" check column pair
# (O = original current-control score at upload time, P = O-M difference)
# to the gradebook sheet, removes the old ad-hoc "очень хочет 4" / "хочу 4"
# notes, and selects the new P4:P32 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 2): relabel / add the new headers.
#    L2 keeps "Сумм" (unchanged text, just shifts shared-string index
#    automatically once the old "очень хочет 4"/"хочу 4" strings are
#    garbage-collected below). Add M2 "ТК_оригинал" and O2
#    "на момент выгрузки в элжуре".
# ---------------------------------------------------------------------
$ws.Range("M2").Value = "ТК_оригинал"
$ws.Range("O2").Value = "на момент выгрузки в элжуре"

# ---------------------------------------------------------------------
# 2. Remove the old one-off remarks that used the now-retired shared
#    strings "очень хочет 4" (N25, N26) and "хочу 4" (O31, overwritten
#    with a numeric value below).
# ---------------------------------------------------------------------
$ws.Range("N25").ClearContents()
$ws.Range("N26").ClearContents()

# ---------------------------------------------------------------------
# 3. New column O ("ТК_оригинал" scores) for every student row, styled
#    like the rest of the graded data columns (copy format from C4,
#    which carries the thick-bordered / centered cell style).
# ---------------------------------------------------------------------
$oValues = @{
    4  = 4
    5  = 5
    6  = 4
    7  = 3
    8  = 4
    9  = 5
    10 = 5
    11 = 5
    12 = 5
    13 = 5
    14 = 5
    15 = 5
    16 = 5
    17 = 3
    18 = 4
    19 = 4
    20 = 5
    21 = 3
    22 = 5
    23 = 3
    24 = 5
    25 = 5
    26 = 3
    27 = 5
    28 = 3
    29 = 5
    30 = 3
    31 = 4
    32 = 5
}

for ($r = 4; $r -le 32; $r++) {
    $ws.Cells.Item($r, 15).Value = $oValues[$r]
}

$ws.Range("C4").Copy()
$ws.Range("O4:O32").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. New column P = O - M (difference between the original TK score and
#    the recorded mark). P4 seeds the shared-formula range P5:P32.
# ---------------------------------------------------------------------
$ws.Range("P4").Formula = "=O4-M4"
$ws.Range("P5:P32").Formula = "=O5-M5"

# ---------------------------------------------------------------------
# 5. Row heights: the whole data block (rows 4-32) now uses the taller,
#    thick-top/bottom-ruled row style already used by rows 9-11/24-26,
#    and the trailing blank row 33 grows a top rule to match.
# ---------------------------------------------------------------------
for ($r = 4; $r -le 32; $r++) {
    $ws.Rows.Item($r).RowHeight = 14
}
$ws.Rows.Item(33).RowHeight = 13

# ---------------------------------------------------------------------
# 6. Selection: the new column is selected top-to-bottom (P4:P32) with
#    the active cell on P4, matching the refreshed view state.
# ---------------------------------------------------------------------
$ws.Range("P4:P32").Select()
